# Add a new price column ("Лом_HMS 1/2 80:20, FOB США Восток, $/т") before the
# existing "Лом_3А, FOB РФ Черное море, $/т" column (old column H), shifting
# every column from H through S one slot to the right (H->I, I->J, ... R->S).
# The old "cap" column (old S) is then removed completely (it lands at T right
# after the insert, so deleting T cancels out the shift for everything past
# it and the sheet ends up the same width as before: A:U). Finally the
# "floor" header at T is renamed to "price_diff", and the price_diff /
# predicted_price (U) values are refreshed to their new numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert a brand-new blank column at H; everything from H.. shifts right.
$ws.Columns("H:H").Insert()

# 2) Remove the old "cap" column, which after the insert above now lives at T.
$ws.Columns("T:T").Delete()

# 3) Header + values for the freshly inserted column H.
$ws.Range("H1").Value = "Лом_HMS 1/2 80:20, FOB США Восток, $/т"

$newH = @(368,342,322,322,338,340,336,328,323,320,313,306,313,340,354,358,356,367,387,395,385,393,404,409,400,422,430,434)
for ($i = 0; $i -lt $newH.Length; $i++) {
    $ws.Cells.Item($i + 2, 8).Value = $newH[$i]
}

# 4) Rename the "floor" header (column T) to "price_diff" and refresh values.
$ws.Range("T1").Value = "price_diff"

$newT = @(1500,-1000,-500,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0)
for ($i = 0; $i -lt $newT.Length; $i++) {
    $ws.Cells.Item($i + 2, 20).Value = $newT[$i]
}

# 5) Refresh the "predicted_price" column (U) values.
$newU = @(47500,42000,41000,42500,45200,44500,43500,41200,41000,41700,42100,38600,38500,39000,39600,39600,40000,41100,43600,43700,44100,43600,45000,44300,45300,46100,48600,50600)
for ($i = 0; $i -lt $newU.Length; $i++) {
    $ws.Cells.Item($i + 2, 21).Value = $newU[$i]
}
